# Updated cryptos list values per source diff (Price + Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.798.19"
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = "'2.369.14"
$ws.Range("E3").Value = '  -3.60%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = "'543.23"
$ws.Range("E5").Value = '  -0.60%  '

$ws.Range("D6").Value = "'140.80"
$ws.Range("E6").Value = '  -2.82%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = "'0.552"
$ws.Range("E8").Value = '  -7.12%  '

$ws.Range("D9").Value = "'2.367.30"
$ws.Range("E9").Value = '  -3.64%  '

$ws.Range("E10").Value = '  -1.57%  '

$ws.Range("E11").Value = '  +0.64%  '

$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("D13").Value = "'0.343"
$ws.Range("E13").Value = '  -2.24%  '

$ws.Range("D14").Value = "'25.39"
$ws.Range("E14").Value = '  -1.88%  '

$ws.Range("D15").Value = "'2.801.54"
$ws.Range("E15").Value = '  -3.37%  '

$ws.Range("E16").Value = '  +0.49%  '

$ws.Range("D17").Value = "'60.824.54"
$ws.Range("E17").Value = '  +0.16%  '

$ws.Range("D18").Value = "'2.370.81"
$ws.Range("E18").Value = '  -3.52%  '

$ws.Range("E19").Value = '  -4.09%  '

$ws.Range("D20").Value = "'4.10"
$ws.Range("E20").Value = '  -1.53%  '

$ws.Range("D21").Value = "'316.93"
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").Value = "'6.68"
$ws.Range("E22").Value = '  -3.51%  '

$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").Value = "'1.80"
$ws.Range("E24").Value = '  +3.90%  '

$ws.Range("D25").Value = "'62.82"
$ws.Range("E25").Value = '  -0.31%  '

$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = '  -0.27%  '

$ws.Range("D27").Value = "'2.487.54"
$ws.Range("E27").Value = '  -3.45%  '

$ws.Range("D28").Value = "'0.0₃0923"
$ws.Range("E28").Value = '  -5.65%  '

$ws.Range("E29").Value = '  +1.58%  '

$ws.Range("D30").Value = "'516.58"
$ws.Range("E30").Value = '  -2.10%  '

$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = '  -3.60%  '

$ws.Range("D32").Value = "'7.96"
$ws.Range("E32").Value = '  -3.44%  '

$ws.Range("E33").Value = '  -2.41%  '

$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = '  -2.69%  '

$ws.Range("E35").Value = '  -0.63%  '

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = '  +0.11%  '

$ws.Range("D37").Value = "'5.44"
$ws.Range("E37").Value = '  -6.66%  '

$ws.Range("E38").Value = '  -4.23%  '

$ws.Range("E39").Value = '  +0.03%  '

$ws.Range("D40").Value = "'18.04"
$ws.Range("E40").Value = '  -1.18%  '

$ws.Range("E41").Value = '  +1.70%  '

$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").Value = "'137.20"
$ws.Range("E43").Value = '  -5.08%  '

$ws.Range("D44").Value = "'40.19"
$ws.Range("E44").Value = '  +0.74%  '

$ws.Range("E45").Value = '  -1.61%  '

$ws.Range("D46").Value = "'139.03"
$ws.Range("E46").Value = '  -5.00%  '

$ws.Range("E47").Value = '  -0.50%  '

$ws.Range("D48").Value = "'20.32"
$ws.Range("E48").Value = '  -2.06%  '

$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = '  -2.67%  '

$ws.Range("D50").Value = "'0.574"
$ws.Range("E50").Value = '  -0.94%  '

$ws.Range("D51").Value = "'0.0914"
$ws.Range("E51").Value = '  -2.77%  '
